$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.868.94"
Set-TextValue $ws.Range("E2") "  -0.80%  "
Set-TextValue $ws.Range("D3") "3.806.73"
Set-TextValue $ws.Range("E3") "  -2.53%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.25%  "
Set-TextValue $ws.Range("D5") "598.89"
Set-TextValue $ws.Range("E5") "  -0.36%  "
Set-TextValue $ws.Range("D6") "168.50"
Set-TextValue $ws.Range("E6") "  -2.55%  "
Set-TextValue $ws.Range("D7") "3.806.77"
Set-TextValue $ws.Range("E7") "  -2.55%  "
Set-TextValue $ws.Range("E8") "  +0.10%  "
Set-TextValue $ws.Range("D9") "0.530"
Set-TextValue $ws.Range("E9") "  -0.41%  "
Set-TextValue $ws.Range("D10") "0.165"
Set-TextValue $ws.Range("E10") "  +0.08%  "
Set-TextValue $ws.Range("D11") "6.50"
Set-TextValue $ws.Range("E11") "  +1.12%  "
Set-TextValue $ws.Range("D12") "0.460"
Set-TextValue $ws.Range("E12") "  +0.11%  "
Set-TextValue $ws.Range("D13") "0.0000272"
Set-TextValue $ws.Range("E13") "  +5.22%  "
Set-TextValue $ws.Range("D14") "36.90"
Set-TextValue $ws.Range("E14") "  -1.01%  "
Set-TextValue $ws.Range("D15") "4.452.96"
Set-TextValue $ws.Range("E15") "  -2.33%  "
Set-TextValue $ws.Range("D16") "3.845.20"
Set-TextValue $ws.Range("E16") "  -1.20%  "
Set-TextValue $ws.Range("D17") "18.86"
Set-TextValue $ws.Range("E17") "  +3.66%  "
Set-TextValue $ws.Range("D18") "67.864.85"
Set-TextValue $ws.Range("E18") "  -0.84%  "
Set-TextValue $ws.Range("D19") "7.33"
Set-TextValue $ws.Range("E19") "  -0.81%  "
Set-TextValue $ws.Range("E20") "  +0.57%  "
Set-TextValue $ws.Range("D21") "10.67"
Set-TextValue $ws.Range("E21") "  -2.75%  "
Set-TextValue $ws.Range("D22") "468.36"
Set-TextValue $ws.Range("E22") "  -0.25%  "
Set-TextValue $ws.Range("D23") "0.731"
Set-TextValue $ws.Range("E23") "  -1.41%  "
Set-TextValue $ws.Range("D24") "0.0000151"
Set-TextValue $ws.Range("E24") "  -5.82%  "
Set-TextValue $ws.Range("D25") "83.66"
Set-TextValue $ws.Range("E25") "  +0.12%  "
Set-TextValue $ws.Range("D26") "2.27"
Set-TextValue $ws.Range("E26") "  +1.16%  "
Set-TextValue $ws.Range("D27") "12.20"
Set-TextValue $ws.Range("E27") "  +0.63%  "
Set-TextValue $ws.Range("D28") "10.32"
Set-TextValue $ws.Range("E28") "  +3.34%  "
Set-TextValue $ws.Range("E29") "  -0.11%  "
Set-TextValue $ws.Range("D30") "2.93"
Set-TextValue $ws.Range("E30") "  -1.02%  "
Set-TextValue $ws.Range("D31") "3.963.22"
Set-TextValue $ws.Range("E31") "  -2.21%  "
Set-TextValue $ws.Range("D32") "7.67"
Set-TextValue $ws.Range("E32") "  -1.23%  "
Set-TextValue $ws.Range("D33") "2.27"
Set-TextValue $ws.Range("E33") "  -2.14%  "
Set-TextValue $ws.Range("D34") "30.61"
Set-TextValue $ws.Range("E34") "  -2.68%  "
Set-TextValue $ws.Range("D35") "9.24"
Set-TextValue $ws.Range("E35") "  -2.28%  "
Set-TextValue $ws.Range("D36") "3.776.19"
Set-TextValue $ws.Range("E36") "  -2.50%  "
Set-TextValue $ws.Range("D37") "3.79"
Set-TextValue $ws.Range("E37") "  +1.76%  "
Set-TextValue $ws.Range("D38") "0.105"
Set-TextValue $ws.Range("E38") "  +0.45%  "
Set-TextValue $ws.Range("D39") "5.94"
Set-TextValue $ws.Range("E39") "  +0.21%  "
Set-TextValue $ws.Range("B40") "Kaspa"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D40") "0.139"
Set-TextValue $ws.Range("E40") "  -1.62%  "
Set-TextValue $ws.Range("B41") "Mantle"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D41") "1.01"
Set-TextValue $ws.Range("E41") "  -2.22%  "
Set-TextValue $ws.Range("D42") "1.00"
Set-TextValue $ws.Range("E42") "  +0.30%  "
Set-TextValue $ws.Range("D43") "0.317"
Set-TextValue $ws.Range("E43") "  +1.08%  "
Set-TextValue $ws.Range("E44") "  -0.02%  "
Set-TextValue $ws.Range("D45") "1.97"
Set-TextValue $ws.Range("E45") "  -0.87%  "
Set-TextValue $ws.Range("D46") "8.75"
Set-TextValue $ws.Range("E46") "  +1.20%  "
Set-TextValue $ws.Range("D47") "409.53"
Set-TextValue $ws.Range("E47") "  -3.81%  "
Set-TextValue $ws.Range("D48") "46.43"
Set-TextValue $ws.Range("E48") "  -1.78%  "
Set-TextValue $ws.Range("D49") "0.000280"
Set-TextValue $ws.Range("E49") "  -7.84%  "
Set-TextValue $ws.Range("D50") "142.43"
Set-TextValue $ws.Range("E50") "  -0.92%  "
Set-TextValue $ws.Range("D51") "0.0359"
Set-TextValue $ws.Range("E51") "  +0.16%  "
